# Weekly update: add a new week of "Tomate" price data (fecha 44776) for
# "Comercializadora del Agro de Limarí" ahead of the existing history.
#
# The new observation occupies 3 rows (Primera / Segunda / Tercera quality
# grades) and is inserted right above the existing row 602, which pushes
# every subsequent row down by three (602->605, ..., 689->692) while leaving
# rows 1-601 untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 602:604; Excel shifts rows 602-689 down to 605-692
# and inherits row formatting (incl. the date style on column D) from the
# row immediately above the insertion point.
$ws.Rows("602:604").Insert()

# Row 602 - Primera
$ws.Range("A602").Value = 2
$ws.Range("B602").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C602").Value = "Coquimbo"
$ws.Range("D602").Value = 44776
$ws.Range("E602").Value = 4
$ws.Range("F602").Value = 100112020
$ws.Range("G602").Value = "Tomate"
$ws.Range("H602").Value = "Larga vida"
$ws.Range("I602").Value = "Primera"
$ws.Range("J602").Value = 2400
$ws.Range("K602").Value = 7000
$ws.Range("L602").Value = 8000
$ws.Range("M602").Value = 7500
$ws.Range("N602").Value = "$/bandeja 18 kilos"
$ws.Range("O602").Value = "Provincia de Limarí"
$ws.Range("P602").Value = 417
$ws.Range("Q602").Value = 18
$ws.Range("R602").Value = "Hortaliza"

# Row 603 - Segunda
$ws.Range("A603").Value = 2
$ws.Range("B603").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C603").Value = "Coquimbo"
$ws.Range("D603").Value = 44776
$ws.Range("E603").Value = 4
$ws.Range("F603").Value = 100112020
$ws.Range("G603").Value = "Tomate"
$ws.Range("H603").Value = "Larga vida"
$ws.Range("I603").Value = "Segunda"
$ws.Range("J603").Value = 1800
$ws.Range("K603").Value = 5000
$ws.Range("L603").Value = 6000
$ws.Range("M603").Value = 5500
$ws.Range("N603").Value = "$/bandeja 18 kilos"
$ws.Range("O603").Value = "Provincia de Limarí"
$ws.Range("P603").Value = 306
$ws.Range("Q603").Value = 18
$ws.Range("R603").Value = "Hortaliza"

# Row 604 - Tercera
$ws.Range("A604").Value = 2
$ws.Range("B604").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C604").Value = "Coquimbo"
$ws.Range("D604").Value = 44776
$ws.Range("E604").Value = 4
$ws.Range("F604").Value = 100112020
$ws.Range("G604").Value = "Tomate"
$ws.Range("H604").Value = "Larga vida"
$ws.Range("I604").Value = "Tercera"
$ws.Range("J604").Value = 1700
$ws.Range("K604").Value = 3000
$ws.Range("L604").Value = 4000
$ws.Range("M604").Value = 3500
$ws.Range("N604").Value = "$/bandeja 18 kilos"
$ws.Range("O604").Value = "Provincia de Limarí"
$ws.Range("P604").Value = 194
$ws.Range("Q604").Value = 18
$ws.Range("R604").Value = "Hortaliza"
